# Auto-generated PowerShell Excel COM-interop script
# Applies the column-reorder/rename + label-rename/row-reorder edits
# described by the diff, for each worksheet in the workbook.
#
# Values that look numeric/percent/currency ('444', '61.49%', '$455,735,134')
# must be forced to Text so Excel doesn't silently convert them to numbers,
# matching the original inlineStr (text) cells in the workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Overall ----
$ws = $wb.Worksheets.Item('Overall')
$ws.Cells.Item(1, 1).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 2).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 3).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 4).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).NumberFormat = '@'
$ws.Cells.Item(2, 1).Value = '61.49%'
$ws.Cells.Item(2, 1).Style = 'Normal'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '444'
$ws.Cells.Item(2, 2).Style = 'Normal'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '$455,735,134'
$ws.Cells.Item(2, 3).Style = 'Normal'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '9.73%'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '-6.94%'
$ws.Cells.Item(2, 5).Style = 'Normal'

# ---- Sheet: County ----
$ws = $wb.Worksheets.Item('County')
$ws.Cells.Item(1, 1).Value = 'Geography'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).Value = 'United States'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '67.35%'
$ws.Cells.Item(2, 2).Style = 'Normal'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '103,475'
$ws.Cells.Item(2, 3).Style = 'Normal'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$267,700,640,005'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '9.05%'
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '-12.83%'
$ws.Cells.Item(2, 6).Style = 'Normal'
$ws.Cells.Item(3, 1).Value = 'North Dakota'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '61.49%'
$ws.Cells.Item(3, 2).Style = 'Normal'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '444'
$ws.Cells.Item(3, 3).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$455,735,134'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '9.73%'
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '-6.94%'
$ws.Cells.Item(3, 6).Style = 'Normal'
$ws.Cells.Item(4, 1).Value = 'Adams County'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '50.00%'
$ws.Cells.Item(4, 2).Style = 'Normal'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '2'
$ws.Cells.Item(4, 3).Style = 'Normal'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$2,411,226'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '8.95%'
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '0.78%'
$ws.Cells.Item(4, 6).Style = 'Normal'
$ws.Cells.Item(5, 1).Value = 'Barnes County'
$ws.Cells.Item(5, 2).NumberFormat = '@'
$ws.Cells.Item(5, 2).Value = '62.50%'
$ws.Cells.Item(5, 2).Style = 'Normal'
$ws.Cells.Item(5, 3).NumberFormat = '@'
$ws.Cells.Item(5, 3).Value = '8'
$ws.Cells.Item(5, 3).Style = 'Normal'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '$3,949,890'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '7.94%'
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(5, 6).NumberFormat = '@'
$ws.Cells.Item(5, 6).Value = '-6.66%'
$ws.Cells.Item(5, 6).Style = 'Normal'
$ws.Cells.Item(6, 1).Value = 'Benson County'
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = '50.00%'
$ws.Cells.Item(6, 2).Style = 'Normal'
$ws.Cells.Item(6, 3).NumberFormat = '@'
$ws.Cells.Item(6, 3).Value = '2'
$ws.Cells.Item(6, 3).Style = 'Normal'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '$15,238,133'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '33.24%'
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(6, 6).NumberFormat = '@'
$ws.Cells.Item(6, 6).Value = '-30.72%'
$ws.Cells.Item(6, 6).Style = 'Normal'
$ws.Cells.Item(7, 1).Value = 'Billings County'
$ws.Cells.Item(7, 2).NumberFormat = '@'
$ws.Cells.Item(7, 2).Value = '0.00%'
$ws.Cells.Item(7, 2).Style = 'Normal'
$ws.Cells.Item(7, 3).NumberFormat = '@'
$ws.Cells.Item(7, 3).Value = '2'
$ws.Cells.Item(7, 3).Style = 'Normal'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '$5,500,176'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '50.62%'
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(7, 6).NumberFormat = '@'
$ws.Cells.Item(7, 6).Value = '35.34%'
$ws.Cells.Item(7, 6).Style = 'Normal'
$ws.Cells.Item(8, 1).Value = 'Bottineau County'
$ws.Cells.Item(8, 2).NumberFormat = '@'
$ws.Cells.Item(8, 2).Value = '100.00%'
$ws.Cells.Item(8, 2).Style = 'Normal'
$ws.Cells.Item(8, 3).NumberFormat = '@'
$ws.Cells.Item(8, 3).Value = '3'
$ws.Cells.Item(8, 3).Style = 'Normal'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '$356,895'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '-21.80%'
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(8, 6).NumberFormat = '@'
$ws.Cells.Item(8, 6).Value = '-75.65%'
$ws.Cells.Item(8, 6).Style = 'Normal'
$ws.Cells.Item(9, 1).Value = 'Bowman County'
$ws.Cells.Item(9, 2).NumberFormat = '@'
$ws.Cells.Item(9, 2).Value = '83.33%'
$ws.Cells.Item(9, 2).Style = 'Normal'
$ws.Cells.Item(9, 3).NumberFormat = '@'
$ws.Cells.Item(9, 3).Value = '6'
$ws.Cells.Item(9, 3).Style = 'Normal'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '$4,037,759'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '2.63%'
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(9, 6).NumberFormat = '@'
$ws.Cells.Item(9, 6).Value = '-14.22%'
$ws.Cells.Item(9, 6).Style = 'Normal'
$ws.Cells.Item(10, 1).Value = 'Burleigh County'
$ws.Cells.Item(10, 2).NumberFormat = '@'
$ws.Cells.Item(10, 2).Value = '67.27%'
$ws.Cells.Item(10, 2).Style = 'Normal'
$ws.Cells.Item(10, 3).NumberFormat = '@'
$ws.Cells.Item(10, 3).Value = '55'
$ws.Cells.Item(10, 3).Style = 'Normal'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '$88,646,357'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '12.40%'
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(10, 6).NumberFormat = '@'
$ws.Cells.Item(10, 6).Value = '-8.97%'
$ws.Cells.Item(10, 6).Style = 'Normal'
$ws.Cells.Item(11, 1).Value = 'Cass County'
$ws.Cells.Item(11, 2).NumberFormat = '@'
$ws.Cells.Item(11, 2).Value = '54.43%'
$ws.Cells.Item(11, 2).Style = 'Normal'
$ws.Cells.Item(11, 3).NumberFormat = '@'
$ws.Cells.Item(11, 3).Value = '79'
$ws.Cells.Item(11, 3).Style = 'Normal'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '$79,896,589'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '10.05%'
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(11, 6).NumberFormat = '@'
$ws.Cells.Item(11, 6).Value = '-2.86%'
$ws.Cells.Item(11, 6).Style = 'Normal'
$ws.Cells.Item(12, 1).Value = 'Cavalier County'
$ws.Cells.Item(12, 2).NumberFormat = '@'
$ws.Cells.Item(12, 2).Value = '80.00%'
$ws.Cells.Item(12, 2).Style = 'Normal'
$ws.Cells.Item(12, 3).NumberFormat = '@'
$ws.Cells.Item(12, 3).Value = '5'
$ws.Cells.Item(12, 3).Style = 'Normal'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '$1,573,244'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '6.24%'
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(12, 6).NumberFormat = '@'
$ws.Cells.Item(12, 6).Value = '-3.11%'
$ws.Cells.Item(12, 6).Style = 'Normal'
$ws.Cells.Item(13, 1).Value = 'Dickey County'
$ws.Cells.Item(13, 2).NumberFormat = '@'
$ws.Cells.Item(13, 2).Value = '60.00%'
$ws.Cells.Item(13, 2).Style = 'Normal'
$ws.Cells.Item(13, 3).NumberFormat = '@'
$ws.Cells.Item(13, 3).Value = '5'
$ws.Cells.Item(13, 3).Style = 'Normal'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '$2,041,437'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '8.58%'
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(13, 6).NumberFormat = '@'
$ws.Cells.Item(13, 6).Value = '-12.23%'
$ws.Cells.Item(13, 6).Style = 'Normal'
$ws.Cells.Item(14, 1).Value = 'Divide County'
$ws.Cells.Item(14, 2).NumberFormat = '@'
$ws.Cells.Item(14, 2).Value = '100.00%'
$ws.Cells.Item(14, 2).Style = 'Normal'
$ws.Cells.Item(14, 3).NumberFormat = '@'
$ws.Cells.Item(14, 3).Value = '1'
$ws.Cells.Item(14, 3).Style = 'Normal'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '$771,047'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '-12.43%'
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(14, 6).NumberFormat = '@'
$ws.Cells.Item(14, 6).Value = '-21.25%'
$ws.Cells.Item(14, 6).Style = 'Normal'
$ws.Cells.Item(15, 1).Value = 'Dunn County'
$ws.Cells.Item(15, 2).NumberFormat = '@'
$ws.Cells.Item(15, 2).Value = '50.00%'
$ws.Cells.Item(15, 2).Style = 'Normal'
$ws.Cells.Item(15, 3).NumberFormat = '@'
$ws.Cells.Item(15, 3).Value = '4'
$ws.Cells.Item(15, 3).Style = 'Normal'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '$437,622'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '30.16%'
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(15, 6).NumberFormat = '@'
$ws.Cells.Item(15, 6).Value = '6.84%'
$ws.Cells.Item(15, 6).Style = 'Normal'
$ws.Cells.Item(16, 1).Value = 'Eddy County'
$ws.Cells.Item(16, 2).NumberFormat = '@'
$ws.Cells.Item(16, 2).Value = '66.67%'
$ws.Cells.Item(16, 2).Style = 'Normal'
$ws.Cells.Item(16, 3).NumberFormat = '@'
$ws.Cells.Item(16, 3).Value = '3'
$ws.Cells.Item(16, 3).Style = 'Normal'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '$1,179,313'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '6.92%'
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(16, 6).NumberFormat = '@'
$ws.Cells.Item(16, 6).Value = '-7.08%'
$ws.Cells.Item(16, 6).Style = 'Normal'
$ws.Cells.Item(17, 1).Value = 'Emmons County'
$ws.Cells.Item(17, 2).NumberFormat = '@'
$ws.Cells.Item(17, 2).Value = '100.00%'
$ws.Cells.Item(17, 2).Style = 'Normal'
$ws.Cells.Item(17, 3).NumberFormat = '@'
$ws.Cells.Item(17, 3).Value = '3'
$ws.Cells.Item(17, 3).Style = 'Normal'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '$2,421,293'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '10.03%'
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(17, 6).NumberFormat = '@'
$ws.Cells.Item(17, 6).Value = '-7.50%'
$ws.Cells.Item(17, 6).Style = 'Normal'
$ws.Cells.Item(18, 1).Value = 'Foster County'
$ws.Cells.Item(18, 2).NumberFormat = '@'
$ws.Cells.Item(18, 2).Value = '33.33%'
$ws.Cells.Item(18, 2).Style = 'Normal'
$ws.Cells.Item(18, 3).NumberFormat = '@'
$ws.Cells.Item(18, 3).Value = '3'
$ws.Cells.Item(18, 3).Style = 'Normal'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '$764,669'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '9.92%'
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(18, 6).NumberFormat = '@'
$ws.Cells.Item(18, 6).Value = '5.23%'
$ws.Cells.Item(18, 6).Style = 'Normal'
$ws.Cells.Item(19, 1).Value = 'Golden Valley County'
$ws.Cells.Item(19, 2).NumberFormat = '@'
$ws.Cells.Item(19, 2).Value = '100.00%'
$ws.Cells.Item(19, 2).Style = 'Normal'
$ws.Cells.Item(19, 3).NumberFormat = '@'
$ws.Cells.Item(19, 3).Value = '3'
$ws.Cells.Item(19, 3).Style = 'Normal'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '$274,985'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '2.36%'
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(19, 6).NumberFormat = '@'
$ws.Cells.Item(19, 6).Value = '-13.33%'
$ws.Cells.Item(19, 6).Style = 'Normal'
$ws.Cells.Item(20, 1).Value = 'Grand Forks County'
$ws.Cells.Item(20, 2).NumberFormat = '@'
$ws.Cells.Item(20, 2).Value = '62.22%'
$ws.Cells.Item(20, 2).Style = 'Normal'
$ws.Cells.Item(20, 3).NumberFormat = '@'
$ws.Cells.Item(20, 3).Value = '45'
$ws.Cells.Item(20, 3).Style = 'Normal'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '$32,815,451'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '6.64%'
$ws.Cells.Item(20, 5).Style = 'Normal'
$ws.Cells.Item(20, 6).NumberFormat = '@'
$ws.Cells.Item(20, 6).Value = '-12.23%'
$ws.Cells.Item(20, 6).Style = 'Normal'
$ws.Cells.Item(21, 1).Value = 'Grant County'
$ws.Cells.Item(21, 2).NumberFormat = '@'
$ws.Cells.Item(21, 2).Value = '60.00%'
$ws.Cells.Item(21, 2).Style = 'Normal'
$ws.Cells.Item(21, 3).NumberFormat = '@'
$ws.Cells.Item(21, 3).Value = '5'
$ws.Cells.Item(21, 3).Style = 'Normal'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '$1,364,485'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '5.18%'
$ws.Cells.Item(21, 5).Style = 'Normal'
$ws.Cells.Item(21, 6).NumberFormat = '@'
$ws.Cells.Item(21, 6).Value = '-8.74%'
$ws.Cells.Item(21, 6).Style = 'Normal'
$ws.Cells.Item(22, 1).Value = 'Griggs County'
$ws.Cells.Item(22, 2).NumberFormat = '@'
$ws.Cells.Item(22, 2).Value = '33.33%'
$ws.Cells.Item(22, 2).Style = 'Normal'
$ws.Cells.Item(22, 3).NumberFormat = '@'
$ws.Cells.Item(22, 3).Value = '3'
$ws.Cells.Item(22, 3).Style = 'Normal'
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '$1,198,822'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '14.95%'
$ws.Cells.Item(22, 5).Style = 'Normal'
$ws.Cells.Item(22, 6).NumberFormat = '@'
$ws.Cells.Item(22, 6).Value = '3.64%'
$ws.Cells.Item(22, 6).Style = 'Normal'
$ws.Cells.Item(23, 1).Value = 'LaMoure County'
$ws.Cells.Item(23, 2).NumberFormat = '@'
$ws.Cells.Item(23, 2).Value = '80.00%'
$ws.Cells.Item(23, 2).Style = 'Normal'
$ws.Cells.Item(23, 3).NumberFormat = '@'
$ws.Cells.Item(23, 3).Value = '5'
$ws.Cells.Item(23, 3).Style = 'Normal'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '$899,601'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '14.25%'
$ws.Cells.Item(23, 5).Style = 'Normal'
$ws.Cells.Item(23, 6).NumberFormat = '@'
$ws.Cells.Item(23, 6).Value = '-3.72%'
$ws.Cells.Item(23, 6).Style = 'Normal'
$ws.Cells.Item(24, 1).Value = 'Logan County'
$ws.Cells.Item(24, 2).NumberFormat = '@'
$ws.Cells.Item(24, 2).Value = '66.67%'
$ws.Cells.Item(24, 2).Style = 'Normal'
$ws.Cells.Item(24, 3).NumberFormat = '@'
$ws.Cells.Item(24, 3).Value = '3'
$ws.Cells.Item(24, 3).Style = 'Normal'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '$201,123'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '1.21%'
$ws.Cells.Item(24, 5).Style = 'Normal'
$ws.Cells.Item(24, 6).NumberFormat = '@'
$ws.Cells.Item(24, 6).Value = '-1.95%'
$ws.Cells.Item(24, 6).Style = 'Normal'
$ws.Cells.Item(25, 1).Value = 'McHenry County'
$ws.Cells.Item(25, 2).NumberFormat = '@'
$ws.Cells.Item(25, 2).Value = '50.00%'
$ws.Cells.Item(25, 2).Style = 'Normal'
$ws.Cells.Item(25, 3).NumberFormat = '@'
$ws.Cells.Item(25, 3).Value = '2'
$ws.Cells.Item(25, 3).Style = 'Normal'
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '$334,836'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '27.83%'
$ws.Cells.Item(25, 5).Style = 'Normal'
$ws.Cells.Item(25, 6).NumberFormat = '@'
$ws.Cells.Item(25, 6).Value = '4.97%'
$ws.Cells.Item(25, 6).Style = 'Normal'
$ws.Cells.Item(26, 1).Value = 'McIntosh County'
$ws.Cells.Item(26, 2).NumberFormat = '@'
$ws.Cells.Item(26, 2).Value = '100.00%'
$ws.Cells.Item(26, 2).Style = 'Normal'
$ws.Cells.Item(26, 3).NumberFormat = '@'
$ws.Cells.Item(26, 3).Value = '3'
$ws.Cells.Item(26, 3).Style = 'Normal'
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '$3,465,675'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '8.71%'
$ws.Cells.Item(26, 5).Style = 'Normal'
$ws.Cells.Item(26, 6).NumberFormat = '@'
$ws.Cells.Item(26, 6).Value = '-6.20%'
$ws.Cells.Item(26, 6).Style = 'Normal'
$ws.Cells.Item(27, 1).Value = 'McKenzie County'
$ws.Cells.Item(27, 2).NumberFormat = '@'
$ws.Cells.Item(27, 2).Value = '50.00%'
$ws.Cells.Item(27, 2).Style = 'Normal'
$ws.Cells.Item(27, 3).NumberFormat = '@'
$ws.Cells.Item(27, 3).Value = '8'
$ws.Cells.Item(27, 3).Style = 'Normal'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '$6,049,620'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '28.77%'
$ws.Cells.Item(27, 5).Style = 'Normal'
$ws.Cells.Item(27, 6).NumberFormat = '@'
$ws.Cells.Item(27, 6).Value = '0.79%'
$ws.Cells.Item(27, 6).Style = 'Normal'
$ws.Cells.Item(28, 1).Value = 'McLean County'
$ws.Cells.Item(28, 2).NumberFormat = '@'
$ws.Cells.Item(28, 2).Value = '100.00%'
$ws.Cells.Item(28, 2).Style = 'Normal'
$ws.Cells.Item(28, 3).NumberFormat = '@'
$ws.Cells.Item(28, 3).Value = '5'
$ws.Cells.Item(28, 3).Style = 'Normal'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '$5,214,120'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '3.23%'
$ws.Cells.Item(28, 5).Style = 'Normal'
$ws.Cells.Item(28, 6).NumberFormat = '@'
$ws.Cells.Item(28, 6).Value = '-42.47%'
$ws.Cells.Item(28, 6).Style = 'Normal'
$ws.Cells.Item(29, 1).Value = 'Mercer County'
$ws.Cells.Item(29, 2).NumberFormat = '@'
$ws.Cells.Item(29, 2).Value = '62.50%'
$ws.Cells.Item(29, 2).Style = 'Normal'
$ws.Cells.Item(29, 3).NumberFormat = '@'
$ws.Cells.Item(29, 3).Value = '8'
$ws.Cells.Item(29, 3).Style = 'Normal'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '$7,104,051'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '6.70%'
$ws.Cells.Item(29, 5).Style = 'Normal'
$ws.Cells.Item(29, 6).NumberFormat = '@'
$ws.Cells.Item(29, 6).Value = '-12.17%'
$ws.Cells.Item(29, 6).Style = 'Normal'
$ws.Cells.Item(30, 1).Value = 'Morton County'
$ws.Cells.Item(30, 2).NumberFormat = '@'
$ws.Cells.Item(30, 2).Value = '64.29%'
$ws.Cells.Item(30, 2).Style = 'Normal'
$ws.Cells.Item(30, 3).NumberFormat = '@'
$ws.Cells.Item(30, 3).Value = '14'
$ws.Cells.Item(30, 3).Style = 'Normal'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '$5,836,226'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '8.37%'
$ws.Cells.Item(30, 5).Style = 'Normal'
$ws.Cells.Item(30, 6).NumberFormat = '@'
$ws.Cells.Item(30, 6).Value = '-4.36%'
$ws.Cells.Item(30, 6).Style = 'Normal'
$ws.Cells.Item(31, 1).Value = 'Mountrail County'
$ws.Cells.Item(31, 2).NumberFormat = '@'
$ws.Cells.Item(31, 2).Value = '25.00%'
$ws.Cells.Item(31, 2).Style = 'Normal'
$ws.Cells.Item(31, 3).NumberFormat = '@'
$ws.Cells.Item(31, 3).Value = '4'
$ws.Cells.Item(31, 3).Style = 'Normal'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '$499,183'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '12.03%'
$ws.Cells.Item(31, 5).Style = 'Normal'
$ws.Cells.Item(31, 6).NumberFormat = '@'
$ws.Cells.Item(31, 6).Value = '7.98%'
$ws.Cells.Item(31, 6).Style = 'Normal'
$ws.Cells.Item(32, 1).Value = 'Nelson County'
$ws.Cells.Item(32, 2).NumberFormat = '@'
$ws.Cells.Item(32, 2).Value = '60.00%'
$ws.Cells.Item(32, 2).Style = 'Normal'
$ws.Cells.Item(32, 3).NumberFormat = '@'
$ws.Cells.Item(32, 3).Value = '5'
$ws.Cells.Item(32, 3).Style = 'Normal'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '$1,084,101'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '4.05%'
$ws.Cells.Item(32, 5).Style = 'Normal'
$ws.Cells.Item(32, 6).NumberFormat = '@'
$ws.Cells.Item(32, 6).Value = '-5.31%'
$ws.Cells.Item(32, 6).Style = 'Normal'
$ws.Cells.Item(33, 1).Value = 'Pembina County'
$ws.Cells.Item(33, 2).NumberFormat = '@'
$ws.Cells.Item(33, 2).Value = '75.00%'
$ws.Cells.Item(33, 2).Style = 'Normal'
$ws.Cells.Item(33, 3).NumberFormat = '@'
$ws.Cells.Item(33, 3).Value = '8'
$ws.Cells.Item(33, 3).Style = 'Normal'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '$1,321,827'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '9.35%'
$ws.Cells.Item(33, 5).Style = 'Normal'
$ws.Cells.Item(33, 6).NumberFormat = '@'
$ws.Cells.Item(33, 6).Value = '-7.94%'
$ws.Cells.Item(33, 6).Style = 'Normal'
$ws.Cells.Item(34, 1).Value = 'Pierce County'
$ws.Cells.Item(34, 2).NumberFormat = '@'
$ws.Cells.Item(34, 2).Value = '33.33%'
$ws.Cells.Item(34, 2).Style = 'Normal'
$ws.Cells.Item(34, 3).NumberFormat = '@'
$ws.Cells.Item(34, 3).Value = '6'
$ws.Cells.Item(34, 3).Style = 'Normal'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '$8,261,047'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '13.49%'
$ws.Cells.Item(34, 5).Style = 'Normal'
$ws.Cells.Item(34, 6).NumberFormat = '@'
$ws.Cells.Item(34, 6).Value = '1.07%'
$ws.Cells.Item(34, 6).Style = 'Normal'
$ws.Cells.Item(35, 1).Value = 'Ramsey County'
$ws.Cells.Item(35, 2).NumberFormat = '@'
$ws.Cells.Item(35, 2).Value = '60.00%'
$ws.Cells.Item(35, 2).Style = 'Normal'
$ws.Cells.Item(35, 3).NumberFormat = '@'
$ws.Cells.Item(35, 3).Value = '10'
$ws.Cells.Item(35, 3).Style = 'Normal'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '$5,349,243'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '14.41%'
$ws.Cells.Item(35, 5).Style = 'Normal'
$ws.Cells.Item(35, 6).NumberFormat = '@'
$ws.Cells.Item(35, 6).Value = '-11.97%'
$ws.Cells.Item(35, 6).Style = 'Normal'
$ws.Cells.Item(36, 1).Value = 'Ransom County'
$ws.Cells.Item(36, 2).NumberFormat = '@'
$ws.Cells.Item(36, 2).Value = '60.00%'
$ws.Cells.Item(36, 2).Style = 'Normal'
$ws.Cells.Item(36, 3).NumberFormat = '@'
$ws.Cells.Item(36, 3).Value = '5'
$ws.Cells.Item(36, 3).Style = 'Normal'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '$1,477,295'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '9.98%'
$ws.Cells.Item(36, 5).Style = 'Normal'
$ws.Cells.Item(36, 6).NumberFormat = '@'
$ws.Cells.Item(36, 6).Value = '-13.92%'
$ws.Cells.Item(36, 6).Style = 'Normal'
$ws.Cells.Item(37, 1).Value = 'Renville County'
$ws.Cells.Item(37, 2).NumberFormat = '@'
$ws.Cells.Item(37, 2).Value = '100.00%'
$ws.Cells.Item(37, 2).Style = 'Normal'
$ws.Cells.Item(37, 3).NumberFormat = '@'
$ws.Cells.Item(37, 3).Value = '1'
$ws.Cells.Item(37, 3).Style = 'Normal'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '$88,287'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '-9.38%'
$ws.Cells.Item(37, 5).Style = 'Normal'
$ws.Cells.Item(37, 6).NumberFormat = '@'
$ws.Cells.Item(37, 6).Value = '-28.50%'
$ws.Cells.Item(37, 6).Style = 'Normal'
$ws.Cells.Item(38, 1).Value = 'Richland County'
$ws.Cells.Item(38, 2).NumberFormat = '@'
$ws.Cells.Item(38, 2).Value = '60.00%'
$ws.Cells.Item(38, 2).Style = 'Normal'
$ws.Cells.Item(38, 3).NumberFormat = '@'
$ws.Cells.Item(38, 3).Value = '10'
$ws.Cells.Item(38, 3).Style = 'Normal'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '$7,850,121'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '3.95%'
$ws.Cells.Item(38, 5).Style = 'Normal'
$ws.Cells.Item(38, 6).NumberFormat = '@'
$ws.Cells.Item(38, 6).Value = '-9.86%'
$ws.Cells.Item(38, 6).Style = 'Normal'
$ws.Cells.Item(39, 1).Value = 'Rolette County'
$ws.Cells.Item(39, 2).NumberFormat = '@'
$ws.Cells.Item(39, 2).Value = '88.89%'
$ws.Cells.Item(39, 2).Style = 'Normal'
$ws.Cells.Item(39, 3).NumberFormat = '@'
$ws.Cells.Item(39, 3).Value = '9'
$ws.Cells.Item(39, 3).Style = 'Normal'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '$32,176,559'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '4.44%'
$ws.Cells.Item(39, 5).Style = 'Normal'
$ws.Cells.Item(39, 6).NumberFormat = '@'
$ws.Cells.Item(39, 6).Value = '-56.68%'
$ws.Cells.Item(39, 6).Style = 'Normal'
$ws.Cells.Item(40, 1).Value = 'Sargent County'
$ws.Cells.Item(40, 2).NumberFormat = '@'
$ws.Cells.Item(40, 2).Value = '50.00%'
$ws.Cells.Item(40, 2).Style = 'Normal'
$ws.Cells.Item(40, 3).NumberFormat = '@'
$ws.Cells.Item(40, 3).Value = '2'
$ws.Cells.Item(40, 3).Style = 'Normal'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '$689,324'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '19.06%'
$ws.Cells.Item(40, 5).Style = 'Normal'
$ws.Cells.Item(40, 6).NumberFormat = '@'
$ws.Cells.Item(40, 6).Value = '-5.05%'
$ws.Cells.Item(40, 6).Style = 'Normal'
$ws.Cells.Item(41, 1).Value = 'Sioux County'
$ws.Cells.Item(41, 2).NumberFormat = '@'
$ws.Cells.Item(41, 2).Value = '100.00%'
$ws.Cells.Item(41, 2).Style = 'Normal'
$ws.Cells.Item(41, 3).NumberFormat = '@'
$ws.Cells.Item(41, 3).Value = '2'
$ws.Cells.Item(41, 3).Style = 'Normal'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '$16,362,686'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '-6.28%'
$ws.Cells.Item(41, 5).Style = 'Normal'
$ws.Cells.Item(41, 6).NumberFormat = '@'
$ws.Cells.Item(41, 6).Value = '-96.72%'
$ws.Cells.Item(41, 6).Style = 'Normal'
$ws.Cells.Item(42, 1).Value = 'Stark County'
$ws.Cells.Item(42, 2).NumberFormat = '@'
$ws.Cells.Item(42, 2).Value = '71.43%'
$ws.Cells.Item(42, 2).Style = 'Normal'
$ws.Cells.Item(42, 3).NumberFormat = '@'
$ws.Cells.Item(42, 3).Value = '21'
$ws.Cells.Item(42, 3).Style = 'Normal'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '$24,095,380'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '12.98%'
$ws.Cells.Item(42, 5).Style = 'Normal'
$ws.Cells.Item(42, 6).NumberFormat = '@'
$ws.Cells.Item(42, 6).Value = '-23.56%'
$ws.Cells.Item(42, 6).Style = 'Normal'
$ws.Cells.Item(43, 1).Value = 'Steele County'
$ws.Cells.Item(43, 2).NumberFormat = '@'
$ws.Cells.Item(43, 2).Value = '0.00%'
$ws.Cells.Item(43, 2).Style = 'Normal'
$ws.Cells.Item(43, 3).NumberFormat = '@'
$ws.Cells.Item(43, 3).Value = '1'
$ws.Cells.Item(43, 3).Style = 'Normal'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '$3,000'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '5.77%'
$ws.Cells.Item(43, 5).Style = 'Normal'
$ws.Cells.Item(43, 6).NumberFormat = '@'
$ws.Cells.Item(43, 6).Value = '4.50%'
$ws.Cells.Item(43, 6).Style = 'Normal'
$ws.Cells.Item(44, 1).Value = 'Stutsman County'
$ws.Cells.Item(44, 2).NumberFormat = '@'
$ws.Cells.Item(44, 2).Value = '50.00%'
$ws.Cells.Item(44, 2).Style = 'Normal'
$ws.Cells.Item(44, 3).NumberFormat = '@'
$ws.Cells.Item(44, 3).Value = '12'
$ws.Cells.Item(44, 3).Style = 'Normal'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '$10,499,456'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '18.34%'
$ws.Cells.Item(44, 5).Style = 'Normal'
$ws.Cells.Item(44, 6).NumberFormat = '@'
$ws.Cells.Item(44, 6).Value = '-0.20%'
$ws.Cells.Item(44, 6).Style = 'Normal'
$ws.Cells.Item(45, 1).Value = 'Towner County'
$ws.Cells.Item(45, 2).NumberFormat = '@'
$ws.Cells.Item(45, 2).Value = '50.00%'
$ws.Cells.Item(45, 2).Style = 'Normal'
$ws.Cells.Item(45, 3).NumberFormat = '@'
$ws.Cells.Item(45, 3).Value = '2'
$ws.Cells.Item(45, 3).Style = 'Normal'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '$3,058,574'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '48.65%'
$ws.Cells.Item(45, 5).Style = 'Normal'
$ws.Cells.Item(45, 6).NumberFormat = '@'
$ws.Cells.Item(45, 6).Value = '37.99%'
$ws.Cells.Item(45, 6).Style = 'Normal'
$ws.Cells.Item(46, 1).Value = 'Traill County'
$ws.Cells.Item(46, 2).NumberFormat = '@'
$ws.Cells.Item(46, 2).Value = '60.00%'
$ws.Cells.Item(46, 2).Style = 'Normal'
$ws.Cells.Item(46, 3).NumberFormat = '@'
$ws.Cells.Item(46, 3).Value = '5'
$ws.Cells.Item(46, 3).Style = 'Normal'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '$1,125,092'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '13.02%'
$ws.Cells.Item(46, 5).Style = 'Normal'
$ws.Cells.Item(46, 6).NumberFormat = '@'
$ws.Cells.Item(46, 6).Value = '-2.29%'
$ws.Cells.Item(46, 6).Style = 'Normal'
$ws.Cells.Item(47, 1).Value = 'Walsh County'
$ws.Cells.Item(47, 2).NumberFormat = '@'
$ws.Cells.Item(47, 2).Value = '54.55%'
$ws.Cells.Item(47, 2).Style = 'Normal'
$ws.Cells.Item(47, 3).NumberFormat = '@'
$ws.Cells.Item(47, 3).Value = '11'
$ws.Cells.Item(47, 3).Style = 'Normal'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '$8,603,511'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '9.41%'
$ws.Cells.Item(47, 5).Style = 'Normal'
$ws.Cells.Item(47, 6).NumberFormat = '@'
$ws.Cells.Item(47, 6).Value = '-1.69%'
$ws.Cells.Item(47, 6).Style = 'Normal'
$ws.Cells.Item(48, 1).Value = 'Ward County'
$ws.Cells.Item(48, 2).NumberFormat = '@'
$ws.Cells.Item(48, 2).Value = '53.12%'
$ws.Cells.Item(48, 2).Style = 'Normal'
$ws.Cells.Item(48, 3).NumberFormat = '@'
$ws.Cells.Item(48, 3).Value = '32'
$ws.Cells.Item(48, 3).Style = 'Normal'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '$53,010,662'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '11.28%'
$ws.Cells.Item(48, 5).Style = 'Normal'
$ws.Cells.Item(48, 6).NumberFormat = '@'
$ws.Cells.Item(48, 6).Value = '-7.51%'
$ws.Cells.Item(48, 6).Style = 'Normal'
$ws.Cells.Item(49, 1).Value = 'Wells County'
$ws.Cells.Item(49, 2).NumberFormat = '@'
$ws.Cells.Item(49, 2).Value = '100.00%'
$ws.Cells.Item(49, 2).Style = 'Normal'
$ws.Cells.Item(49, 3).NumberFormat = '@'
$ws.Cells.Item(49, 3).Value = '1'
$ws.Cells.Item(49, 3).Style = 'Normal'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '$181,000'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '12.07%'
$ws.Cells.Item(49, 5).Style = 'Normal'
$ws.Cells.Item(49, 6).NumberFormat = '@'
$ws.Cells.Item(49, 6).Value = '-18.51%'
$ws.Cells.Item(49, 6).Style = 'Normal'
$ws.Cells.Item(50, 1).Value = 'Williams County'
$ws.Cells.Item(50, 2).NumberFormat = '@'
$ws.Cells.Item(50, 2).Value = '50.00%'
$ws.Cells.Item(50, 2).Style = 'Normal'
$ws.Cells.Item(50, 3).NumberFormat = '@'
$ws.Cells.Item(50, 3).Value = '12'
$ws.Cells.Item(50, 3).Style = 'Normal'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '$6,014,141'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '13.94%'
$ws.Cells.Item(50, 5).Style = 'Normal'
$ws.Cells.Item(50, 6).NumberFormat = '@'
$ws.Cells.Item(50, 6).Value = '-1.15%'
$ws.Cells.Item(50, 6).Style = 'Normal'

# ---- Sheet: Congressional District ----
$ws = $wb.Worksheets.Item('Congressional District')
$ws.Cells.Item(1, 1).Value = 'Geography'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).Value = 'United States'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '67.35%'
$ws.Cells.Item(2, 2).Style = 'Normal'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '103,475'
$ws.Cells.Item(2, 3).Style = 'Normal'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$267,700,640,005'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '9.05%'
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '-12.83%'
$ws.Cells.Item(2, 6).Style = 'Normal'
$ws.Cells.Item(3, 1).Value = 'North Dakota'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '61.49%'
$ws.Cells.Item(3, 2).Style = 'Normal'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '444'
$ws.Cells.Item(3, 3).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$455,735,134'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '9.73%'
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '-6.94%'
$ws.Cells.Item(3, 6).Style = 'Normal'
$ws.Cells.Item(4, 1).Value = 'Congressional District (at Large)'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '61.49%'
$ws.Cells.Item(4, 2).Style = 'Normal'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '444'
$ws.Cells.Item(4, 3).Style = 'Normal'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$455,735,134'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '9.73%'
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '-6.94%'
$ws.Cells.Item(4, 6).Style = 'Normal'

# ---- Sheet: Size ----
$ws = $wb.Worksheets.Item('Size')
$ws.Cells.Item(1, 1).Value = 'Size'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).Value = 'Between $100K and $499K'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '64.33%'
$ws.Cells.Item(2, 2).Style = 'Normal'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '157'
$ws.Cells.Item(2, 3).Style = 'Normal'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$15,389,793'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '9.98%'
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '-11.99%'
$ws.Cells.Item(2, 6).Style = 'Normal'
$ws.Cells.Item(3, 1).Value = 'Between $1M and $4.99M'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '69.00%'
$ws.Cells.Item(3, 2).Style = 'Normal'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '100'
$ws.Cells.Item(3, 3).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$81,726,133'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '7.71%'
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '-9.23%'
$ws.Cells.Item(3, 6).Style = 'Normal'
$ws.Cells.Item(4, 1).Value = 'Between $500K and $999K'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '68.12%'
$ws.Cells.Item(4, 2).Style = 'Normal'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '69'
$ws.Cells.Item(4, 3).Style = 'Normal'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$31,348,865'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '16.52%'
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '-8.25%'
$ws.Cells.Item(4, 6).Style = 'Normal'
$ws.Cells.Item(5, 1).Value = 'Between $5M and $9.99M'
$ws.Cells.Item(5, 2).NumberFormat = '@'
$ws.Cells.Item(5, 2).Value = '58.62%'
$ws.Cells.Item(5, 2).Style = 'Normal'
$ws.Cells.Item(5, 3).NumberFormat = '@'
$ws.Cells.Item(5, 3).Value = '29'
$ws.Cells.Item(5, 3).Style = 'Normal'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '$50,210,634'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '8.71%'
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(5, 6).NumberFormat = '@'
$ws.Cells.Item(5, 6).Value = '-6.20%'
$ws.Cells.Item(5, 6).Style = 'Normal'
$ws.Cells.Item(6, 1).Value = 'Greater than $10M'
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = '51.56%'
$ws.Cells.Item(6, 2).Style = 'Normal'
$ws.Cells.Item(6, 3).NumberFormat = '@'
$ws.Cells.Item(6, 3).Value = '64'
$ws.Cells.Item(6, 3).Style = 'Normal'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '$275,912,964'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '6.48%'
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(6, 6).NumberFormat = '@'
$ws.Cells.Item(6, 6).Value = '-1.52%'
$ws.Cells.Item(6, 6).Style = 'Normal'
$ws.Cells.Item(7, 1).Value = 'Less than $100K'
$ws.Cells.Item(7, 2).NumberFormat = '@'
$ws.Cells.Item(7, 2).Value = '24.00%'
$ws.Cells.Item(7, 2).Style = 'Normal'
$ws.Cells.Item(7, 3).NumberFormat = '@'
$ws.Cells.Item(7, 3).Value = '25'
$ws.Cells.Item(7, 3).Style = 'Normal'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '$1,146,745'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '52.15%'
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(7, 6).NumberFormat = '@'
$ws.Cells.Item(7, 6).Value = '20.55%'
$ws.Cells.Item(7, 6).Style = 'Normal'
$ws.Cells.Item(8, 1).Value = 'Total'
$ws.Cells.Item(8, 2).NumberFormat = '@'
$ws.Cells.Item(8, 2).Value = '61.49%'
$ws.Cells.Item(8, 2).Style = 'Normal'
$ws.Cells.Item(8, 3).NumberFormat = '@'
$ws.Cells.Item(8, 3).Value = '444'
$ws.Cells.Item(8, 3).Style = 'Normal'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '$455,735,134'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '9.73%'
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(8, 6).NumberFormat = '@'
$ws.Cells.Item(8, 6).Value = '-6.94%'
$ws.Cells.Item(8, 6).Style = 'Normal'

# ---- Sheet: Subsector ----
$ws = $wb.Worksheets.Item('Subsector')
$ws.Cells.Item(1, 1).Value = 'Subsector'
$ws.Cells.Item(1, 2).Value = 'Share of 990 filers with government grants at risk'
$ws.Cells.Item(1, 3).Value = 'Number of 990 filers with government grants'
$ws.Cells.Item(1, 4).Value = 'Total government grants ($)'
$ws.Cells.Item(1, 5).Value = 'Size of operating surplus with government grants'
$ws.Cells.Item(1, 6).Value = 'Size of operating surplus without government grants'
$ws.Cells.Item(2, 1).Value = 'Arts, Culture, and Humanities'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '42.86%'
$ws.Cells.Item(2, 2).Style = 'Normal'
$ws.Cells.Item(2, 3).NumberFormat = '@'
$ws.Cells.Item(2, 3).Value = '21'
$ws.Cells.Item(2, 3).Style = 'Normal'
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '$8,058,301'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '25.12%'
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(2, 6).NumberFormat = '@'
$ws.Cells.Item(2, 6).Value = '3.82%'
$ws.Cells.Item(2, 6).Style = 'Normal'
$ws.Cells.Item(3, 1).Value = 'Education (Excluding Universities)'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '37.93%'
$ws.Cells.Item(3, 2).Style = 'Normal'
$ws.Cells.Item(3, 3).NumberFormat = '@'
$ws.Cells.Item(3, 3).Value = '29'
$ws.Cells.Item(3, 3).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '$73,220,273'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '15.02%'
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(3, 6).NumberFormat = '@'
$ws.Cells.Item(3, 6).Value = '3.18%'
$ws.Cells.Item(3, 6).Style = 'Normal'
$ws.Cells.Item(4, 1).Value = 'Environment and Animals'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '50.00%'
$ws.Cells.Item(4, 2).Style = 'Normal'
$ws.Cells.Item(4, 3).NumberFormat = '@'
$ws.Cells.Item(4, 3).Value = '16'
$ws.Cells.Item(4, 3).Style = 'Normal'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '$3,792,984'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '23.12%'
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(4, 6).NumberFormat = '@'
$ws.Cells.Item(4, 6).Value = '2.39%'
$ws.Cells.Item(4, 6).Style = 'Normal'
$ws.Cells.Item(5, 1).Value = 'Health (Excluding Hospitals)'
$ws.Cells.Item(5, 2).NumberFormat = '@'
$ws.Cells.Item(5, 2).Value = '64.81%'
$ws.Cells.Item(5, 2).Style = 'Normal'
$ws.Cells.Item(5, 3).NumberFormat = '@'
$ws.Cells.Item(5, 3).Value = '54'
$ws.Cells.Item(5, 3).Style = 'Normal'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '$32,531,582'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '8.26%'
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(5, 6).NumberFormat = '@'
$ws.Cells.Item(5, 6).Value = '-6.29%'
$ws.Cells.Item(5, 6).Style = 'Normal'
$ws.Cells.Item(6, 1).Value = 'Hospitals'
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = '57.14%'
$ws.Cells.Item(6, 2).Style = 'Normal'
$ws.Cells.Item(6, 3).NumberFormat = '@'
$ws.Cells.Item(6, 3).Value = '21'
$ws.Cells.Item(6, 3).Style = 'Normal'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '$56,076,989'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '3.95%'
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(6, 6).NumberFormat = '@'
$ws.Cells.Item(6, 6).Value = '-5.46%'
$ws.Cells.Item(6, 6).Style = 'Normal'
$ws.Cells.Item(7, 1).Value = 'Human Services'
$ws.Cells.Item(7, 2).NumberFormat = '@'
$ws.Cells.Item(7, 2).Value = '70.14%'
$ws.Cells.Item(7, 2).Style = 'Normal'
$ws.Cells.Item(7, 3).NumberFormat = '@'
$ws.Cells.Item(7, 3).Value = '144'
$ws.Cells.Item(7, 3).Style = 'Normal'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '$73,364,618'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '7.61%'
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(7, 6).NumberFormat = '@'
$ws.Cells.Item(7, 6).Value = '-12.51%'
$ws.Cells.Item(7, 6).Style = 'Normal'
$ws.Cells.Item(8, 1).Value = 'Public, Societal Benefit'
$ws.Cells.Item(8, 2).NumberFormat = '@'
$ws.Cells.Item(8, 2).Value = '56.25%'
$ws.Cells.Item(8, 2).Style = 'Normal'
$ws.Cells.Item(8, 3).NumberFormat = '@'
$ws.Cells.Item(8, 3).Value = '32'
$ws.Cells.Item(8, 3).Style = 'Normal'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '$24,604,213'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '25.04%'
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(8, 6).NumberFormat = '@'
$ws.Cells.Item(8, 6).Value = '-3.79%'
$ws.Cells.Item(8, 6).Style = 'Normal'
$ws.Cells.Item(9, 1).Value = 'Religion Related'
$ws.Cells.Item(9, 2).NumberFormat = '@'
$ws.Cells.Item(9, 2).Value = '40.00%'
$ws.Cells.Item(9, 2).Style = 'Normal'
$ws.Cells.Item(9, 3).NumberFormat = '@'
$ws.Cells.Item(9, 3).Value = '10'
$ws.Cells.Item(9, 3).Style = 'Normal'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '$2,821,085'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '14.10%'
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(9, 6).NumberFormat = '@'
$ws.Cells.Item(9, 6).Value = '3.82%'
$ws.Cells.Item(9, 6).Style = 'Normal'
$ws.Cells.Item(10, 1).Value = 'Unclassified'
$ws.Cells.Item(10, 2).NumberFormat = '@'
$ws.Cells.Item(10, 2).Value = '64.60%'
$ws.Cells.Item(10, 2).Style = 'Normal'
$ws.Cells.Item(10, 3).NumberFormat = '@'
$ws.Cells.Item(10, 3).Value = '113'
$ws.Cells.Item(10, 3).Style = 'Normal'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '$146,183,475'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '8.71%'
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(10, 6).NumberFormat = '@'
$ws.Cells.Item(10, 6).Value = '-9.33%'
$ws.Cells.Item(10, 6).Style = 'Normal'
$ws.Cells.Item(11, 1).Value = 'Universities'
$ws.Cells.Item(11, 2).NumberFormat = '@'
$ws.Cells.Item(11, 2).Value = '50.00%'
$ws.Cells.Item(11, 2).Style = 'Normal'
$ws.Cells.Item(11, 3).NumberFormat = '@'
$ws.Cells.Item(11, 3).Value = '4'
$ws.Cells.Item(11, 3).Style = 'Normal'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '$35,081,614'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '19.92%'
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(11, 6).NumberFormat = '@'
$ws.Cells.Item(11, 6).Value = '5.37%'
$ws.Cells.Item(11, 6).Style = 'Normal'
$ws.Cells.Item(12, 1).Value = 'Total'
$ws.Cells.Item(12, 2).NumberFormat = '@'
$ws.Cells.Item(12, 2).Value = '61.49%'
$ws.Cells.Item(12, 2).Style = 'Normal'
$ws.Cells.Item(12, 3).NumberFormat = '@'
$ws.Cells.Item(12, 3).Value = '444'
$ws.Cells.Item(12, 3).Style = 'Normal'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '$455,735,134'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '9.73%'
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(12, 6).NumberFormat = '@'
$ws.Cells.Item(12, 6).Value = '-6.94%'
$ws.Cells.Item(12, 6).Style = 'Normal'
